# Fix mislabeled tag rows ("unknown" tag mixup) on the per-language tag-stats
# sheets. Only the text of the "Tag" column cells changes; counts/percentages
# stay exactly as they were. Target text values were derived directly from the
# authoritative OOXML diff (resolved through the *new* shared-string table).

$wb = $excel.ActiveWorkbook

# --- Vietnamese sheet: rows 16-21, Tag columns A, F, K, P ---
$ws = $wb.Worksheets.Item("Vietnamese")

$ws.Range("A16").Value = "INTJ"
$ws.Range("F16").Value = "INTJ"
$ws.Range("K16").Value = "INTJ"
$ws.Range("P16").Value = "INTJ"

$ws.Range("A17").Value = "PRON"
$ws.Range("F17").Value = "PRON"
$ws.Range("K17").Value = "PRON"
$ws.Range("P17").Value = "PRON"

$ws.Range("A18").Value = "SYM"
$ws.Range("F18").Value = "SYM"
$ws.Range("K18").Value = "SYM"
$ws.Range("P18").Value = "SYM"

$ws.Range("A19").Value = "O"
$ws.Range("F19").Value = "O"
$ws.Range("K19").Value = "O"
$ws.Range("P19").Value = "O"

$ws.Range("A20").Value = "ADV"
$ws.Range("F20").Value = "ADV"
$ws.Range("K20").Value = "ADV"
$ws.Range("P20").Value = "ADV"

$ws.Range("A21").Value = "_"
$ws.Range("F21").Value = "_"
$ws.Range("K21").Value = "_"
$ws.Range("P21").Value = "_"

# --- Thai sheet: rows 18-20, Tag column A ---
$ws = $wb.Worksheets.Item("Thai")

$ws.Range("A18").Value = "O"
$ws.Range("A19").Value = "INTJ"
$ws.Range("A20").Value = "X"

# --- Cantonese sheet: row 14 and rows 18-21, Tag column A ---
$ws = $wb.Worksheets.Item("Cantonese")

$ws.Range("A14").Value = "INTJ"
$ws.Range("A18").Value = "SYM"
$ws.Range("A19").Value = "O"
$ws.Range("A20").Value = "X"
$ws.Range("A21").Value = "_"

# --- Isolating sheet: row 18 (col I), row 19 (cols A,E,M), row 21 (cols A,E,I,M) ---
$ws = $wb.Worksheets.Item("Isolating")

$ws.Range("I18").Value = "INTJ"

$ws.Range("A19").Value = "INTJ"
$ws.Range("E19").Value = "INTJ"
$ws.Range("M19").Value = "INTJ"

$ws.Range("A21").Value = "_"
$ws.Range("E21").Value = "_"
$ws.Range("I21").Value = "_"
$ws.Range("M21").Value = "_"
